$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 274, shifting existing rows 274-343 down
$ws.Rows.Item(274).Insert()

# Populate the new row 274 with data
$ws.Cells.Item(274, 1).Value = 10
$ws.Cells.Item(274, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(274, 3).Value = "La Araucanía"
$ws.Cells.Item(274, 4).Value = 44798
$ws.Cells.Item(274, 5).Value = 9
$ws.Cells.Item(274, 6).Value = 100114013
$ws.Cells.Item(274, 7).Value = "Zanahoria"
$ws.Cells.Item(274, 8).Value = "Sin especificar"
$ws.Cells.Item(274, 9).Value = "Primera"
$ws.Cells.Item(274, 10).Value = 150
$ws.Cells.Item(274, 11).Value = 13000
$ws.Cells.Item(274, 12).Value = 13000
$ws.Cells.Item(274, 13).Value = 13000
$ws.Cells.Item(274, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(274, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(274, 16).Value = 520
$ws.Cells.Item(274, 17).Value = 25
$ws.Cells.Item(274, 18).Value = "Hortaliza"
